$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Update Status fields in the zh-cn and de-de tables
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# Update Latest Handoff Datetime fields
$wsZhCn.Range("H2").Value = "2016-08-26 20:38:19"
$wsDeDe.Range("H2").Value = "2016-08-26 20:38:24"

# Update Overview sheet: zh-cn / de-de status columns, and Latest HO Xliff Generate Date
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-26 20:38:24"

# Set the columns' widths to match the widened Status/zh-cn/de-de columns
$wsOverview.Columns.Item(5).ColumnWidth = 16.3333333333333
$wsOverview.Columns.Item(6).ColumnWidth = 16.3333333333333
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3333333333333
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3333333333333
